# Applies the "removed excess columns in dtr summary and added legends to
# per employee report" edit described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) DTR summary: row 10 (04-28-2015 OB row) had a stray 0.5 in I10 that
#    should not be there. Clear the value but keep the cell's formatting.
# ---------------------------------------------------------------------
$ws.Range("I10").ClearContents()

# ---------------------------------------------------------------------
# 2) Add a "Legends:" section under the existing totals block (rows
#    19-22), reusing the row 24-30 area. Row 24 is the heading; rows
#    25-26 / 27-28 / 29-30 are each a merged color-swatch cell (E) next
#    to a merged bold+underlined description cell (F:P).
# ---------------------------------------------------------------------

# Row 24: "Legends:" heading, same big bold-underline title look as A1.
$ws.Range("E24:P24").Merge()
$ws.Range("A1").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$ws.Range("E24").Value = "Legends:"

# Row 25-26: blue swatch + description (request/remark note)
$ws.Range("E25:E26").Merge()
$ws.Range("E25").Interior.Color = 13411113  # FF29A3CC

$ws.Range("F25:P26").Merge()
$ws.Range("F25").Value = "Employee has request(s)/remark(s) for that day." + [char]10 + "*May incur late and/or undertime depending on his or her time-in and time-out."
$ws.Range("F25").Font.Bold = $true
$ws.Range("F25").Font.Underline = $true
$ws.Range("F25").WrapText = $true

# Row 27-28: orange swatch + description (half-day note)
$ws.Range("E27:E28").Merge()
$ws.Range("E27").Interior.Color = 6737407  # FFFFCC66

$ws.Range("F27:P28").Merge()
$ws.Range("F27").Value = "Employee is considered half-day because of his time-in or time-out."
$ws.Range("F27").Font.Bold = $true
$ws.Range("F27").Font.Underline = $true
$ws.Range("F27").WrapText = $true

# Row 29-30: red swatch + description (absent note)
$ws.Range("E29:E30").Merge()
$ws.Range("E29").Interior.Color = 6184671  # FFDF5E5E

$ws.Range("F29:P30").Merge()
$ws.Range("F29").Value = "Employee has no time-in and therefore, considered as absent."
$ws.Range("F29").Font.Bold = $true
$ws.Range("F29").Font.Underline = $true
$ws.Range("F29").WrapText = $true
